# ---------------------------------------------------------------------------
# Add a new "2022-Q4" worksheet (inserted right after "总计" / before
# "2022-Q3") with its fund-holding table, and update the "总计" summary
# sheet so that the new quarter is reflected there too.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Helper: write a value as TEXT (avoids Excel auto-converting numeric-looking
# strings such as "013869" or "0.18" into real numbers), while keeping the
# cell on the default ("Normal") style once done.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet before "2022-Q3"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row - reuse the existing bold/border style already used for the
# other report headers (总计!B1) instead of building a brand new style.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Index column (A2:A3) - reuse style from the summary sheet's index column.
$summary.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1

# Row 2: fund 013869
Set-TextValue $q4.Range("B2") "013869"
$q4.Range("C2").Value = "创金合信物联网主题股票A"
Set-TextValue $q4.Range("D2") "0.18"
Set-TextValue $q4.Range("E2") "87.02"
Set-TextValue $q4.Range("F2") "3.09"
Set-TextValue $q4.Range("G2") "0.0056"
$q4.Range("H2").Value = 9

# Row 3: fund 013870
Set-TextValue $q4.Range("B3") "013870"
$q4.Range("C3").Value = "创金合信物联网主题股票C"
Set-TextValue $q4.Range("D3") "0.12"
Set-TextValue $q4.Range("E3") "87.02"
Set-TextValue $q4.Range("F3") "3.09"
Set-TextValue $q4.Range("G3") "0.0037"
$q4.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a row so the old "2022-Q3" row is kept
#    (shifted down one row) and the now-vacated row 2 receives the new
#    "2022-Q4" figures. The former "2022-Q1" row (now row 4) gets its index
#    bumped from 1 to 2.
# ---------------------------------------------------------------------------
$summary.Rows.Item(3).Insert()

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.07000000000000001

$summary.Range("A4").Value = 2

# Keep the originally-active sheet ("2022-Q1") selected, as in the source
# workbook, rather than whichever sheet our edits last touched.
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Host "2022-Q4 sheet added and 总计 updated"
